$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1832.8286
$ws.Range("I107").Value = 1323.2258
$ws.Range("K107").Value = 1323.2258
$ws.Range("M107").Value = 596.7742000000001

$ws.Range("H125").Value = 4705.231
$ws.Range("I125").Value = 4585.2856
$ws.Range("K125").Value = 41267.5704
$ws.Range("M125").Value = -38807.5704

$ws.Range("H132").Value = 2624.3333
$ws.Range("I132").Value = 2911.077
$ws.Range("J132").Value = 1878.8
$ws.Range("K132").Value = 8733.231
$ws.Range("L132").Value = 5636.4
$ws.Range("M132").Value = -6203.231
$ws.Range("N132").Value = -10696.4

$ws.Range("H137").Value = 14057.743
$ws.Range("I137").Value = 5096.1816
$ws.Range("J137").Value = 29223.46
$ws.Range("K137").Value = 15288.5448
$ws.Range("L137").Value = 87670.38
$ws.Range("M137").Value = -12738.5448
$ws.Range("N137").Value = -92770.38

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3624.2673
$ws.Range("I32").Value = 2365.3215
$ws.Range("K32").Value = 2365.3215
$ws.Range("M32").Value = -2078.3215

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 456483.6
$ws.Range("I86").Value = 770370.0600000001
$ws.Range("K86").Value = 770370.0600000001
$ws.Range("M86").Value = -769247.0600000001

$ws.Range("H89").Value = 456483.6
$ws.Range("I89").Value = 770370.0600000001
$ws.Range("K89").Value = 3851850.3
$ws.Range("M89").Value = -3846234.3

$ws.Range("H105").Value = 2913.425
$ws.Range("I105").Value = 2797.4783
$ws.Range("K105").Value = 2797.4783
$ws.Range("M105").Value = -1050.4783

$ws.Range("H106").Value = 28780
$ws.Range("J106").Value = 28780
$ws.Range("L106").Value = 28780
$ws.Range("N106").Value = -31304

$ws.Range("H134").Value = 7096.7646
$ws.Range("I134").Value = 3727.0833
$ws.Range("J134").Value = 15184
$ws.Range("K134").Value = 11181.2499
$ws.Range("L134").Value = 45552
$ws.Range("M134").Value = -8646.249899999999
$ws.Range("N134").Value = -50622

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 531.8333
$ws.Range("I5").Value = 531.8333
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 531.8333
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -419.8333
$ws.Range("N5").ClearContents()

$ws.Range("H10").Value = 2271.1428
$ws.Range("I10").Value = 724.75
$ws.Range("K10").Value = 724.75
$ws.Range("M10").Value = -585.75

$ws.Range("H11").Value = 1199.75
$ws.Range("J11").Value = 1199.75
$ws.Range("L11").Value = 1199.75
$ws.Range("N11").Value = -1479.75

$ws.Range("H12").Value = 768
$ws.Range("I12").Value = 94.75
$ws.Range("J12").Value = 1665.6666
$ws.Range("K12").Value = 94.75
$ws.Range("L12").Value = 1665.6666
$ws.Range("M12").Value = 75.25
$ws.Range("N12").Value = -2005.6666

$ws.Range("H13").Value = 14222.056
$ws.Range("I13").Value = 14218.6875
$ws.Range("J13").Value = 14249
$ws.Range("K13").Value = 14218.6875
$ws.Range("L13").Value = 14249
$ws.Range("N13").Value = -14527
$ws.Range("M13").Value = -14079.6875

$ws.Range("H31").Value = 316415.56
$ws.Range("I31").Value = 3272.5789
$ws.Range("J31").Value = 774086.0600000001
$ws.Range("K31").Value = 3272.5789
$ws.Range("L31").Value = 774086.0600000001
$ws.Range("M31").Value = -2977.5789
$ws.Range("N31").Value = -774676.0600000001

$ws.Range("H34").Value = 316415.56
$ws.Range("I34").Value = 3272.5789
$ws.Range("J34").Value = 774086.0600000001
$ws.Range("K34").Value = 3272.5789
$ws.Range("L34").Value = 774086.0600000001
$ws.Range("M34").Value = -3070.5789
$ws.Range("N34").Value = -774490.0600000001

$ws.Range("H58").Value = 3714.36
$ws.Range("I58").Value = 2215.5789
$ws.Range("J58").Value = 8460.5
$ws.Range("K58").Value = 2215.5789
$ws.Range("L58").Value = 8460.5
$ws.Range("M58").Value = -2012.5789
$ws.Range("N58").Value = -8866.5

$ws.Range("H94").Value = 3576.5
$ws.Range("I94").Value = 1864.75
$ws.Range("J94").Value = 7000
$ws.Range("K94").Value = 1864.75
$ws.Range("L94").Value = 7000
$ws.Range("M94").Value = -1413.75
$ws.Range("N94").Value = -7902

$ws.Range("H95").Value = 5251.8
$ws.Range("J95").Value = 5390.778
$ws.Range("L95").Value = 5390.778
$ws.Range("N95").Value = -10882.778

$ws.Range("H132").Value = 24694.691
$ws.Range("I132").Value = 16186.907
$ws.Range("J132").Value = 41323.547
$ws.Range("K132").Value = 48560.721
$ws.Range("L132").Value = 123970.641
$ws.Range("M132").Value = -46030.721
$ws.Range("N132").Value = -129030.641

$ws.Range("H134").Value = 3200.3823
$ws.Range("I134").Value = 2738.1155
$ws.Range("J134").Value = 4702.75
$ws.Range("K134").Value = 8214.3465
$ws.Range("L134").Value = 14108.25
$ws.Range("M134").Value = -5679.3465
$ws.Range("N134").Value = -19178.25

$ws.Range("H136").Value = 3714.36
$ws.Range("I136").Value = 2215.5789
$ws.Range("J136").Value = 8460.5
$ws.Range("K136").Value = 6646.736699999999
$ws.Range("L136").Value = 25381.5
$ws.Range("M136").Value = -4096.736699999999
$ws.Range("N136").Value = -30481.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 514.875
$ws.Range("I8").Value = 514.875
$ws.Range("K8").Value = 1544.625
$ws.Range("M8").Value = -1405.625

$ws.Range("H99").Value = 6490.1665
$ws.Range("I99").Value = 3485.25
$ws.Range("K99").Value = 10455.75
$ws.Range("M99").Value = -8209.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 82485
$ws.Range("J100").Value = 82485
$ws.Range("L100").Value = 82485
$ws.Range("N100").Value = -84649

$ws.Range("H102").Value = 10206110
$ws.Range("I102").Value = 2185.1538
$ws.Range("J102").Value = 142857140
$ws.Range("K102").Value = 2185.1538
$ws.Range("L102").Value = 142857140
$ws.Range("M102").Value = -563.1538
$ws.Range("N102").Value = -142860384

$ws.Range("H132").Value = 17279.25
$ws.Range("I132").Value = 12198.048
$ws.Range("J132").Value = 52847.668
$ws.Range("K132").Value = 36594.144
$ws.Range("L132").Value = 158543.004
$ws.Range("M132").Value = -34064.144
$ws.Range("N132").Value = -163603.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4010.6956
$ws.Range("I7").Value = 3738.5454
$ws.Range("J7").Value = 9998
$ws.Range("K7").Value = 3738.5454
$ws.Range("L7").Value = 9998
$ws.Range("M7").Value = -3626.5454
$ws.Range("N7").Value = -10222

$ws.Range("H97").Value = 41359.875
$ws.Range("J97").Value = 41359.875
$ws.Range("L97").Value = 41359.875
$ws.Range("N97").Value = -43341.875

$ws.Range("H126").Value = 4010.6956
$ws.Range("I126").Value = 3738.5454
$ws.Range("J126").Value = 9998
$ws.Range("K126").Value = 11215.6362
$ws.Range("L126").Value = 29994
$ws.Range("M126").Value = -8745.636200000001
$ws.Range("N126").Value = -34934

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4409.683
$ws.Range("I122").Value = 3197.5151
$ws.Range("J122").Value = 9409.875
$ws.Range("K122").Value = 9592.5453
$ws.Range("L122").Value = 28229.625
$ws.Range("M122").Value = -7142.5453
$ws.Range("N122").Value = -33129.625

$ws.Range("H132").Value = 24895.895
$ws.Range("I132").Value = 23134.477
$ws.Range("K132").Value = 69403.431
$ws.Range("M132").Value = -66873.431

$ws.Range("H136").Value = 1721.1526
$ws.Range("I136").Value = 821.6222
$ws.Range("K136").Value = 2464.8666
$ws.Range("M136").Value = 85.13339999999971
